$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: fix connector part numbers (end card to end cap connector)
$ws.Range("E9").Value = "M55-6001242R"
$ws.Range("F9").Value = "12 Position Receptacle Connector 0.050`" (1.27mm) Surface Mount Gold"
$ws.Range("I9").Value = "952-3835-1-ND"

# Row 23: restore literal dist part number (was referencing wrong shared string)
$ws.Range("I23").Value = "139-QCN-27"

# Selection moved to I9
$ws.Range("I9").Select()
